$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data cells in this sheet are plain text (e.g. "1.00", "0.530")
# even though they look numeric. Assigning such strings straight to .Value
# makes Excel coerce them to numbers (losing the literal formatting), and
# forcing a Text NumberFormat instead would stamp a new cell style that is
# not part of this edit. So each cell is written as a `="literal"` text
# formula first; after all of them are in place a single values-only paste
# collapses every formula back down to the plain literal string, matching
# the original inline-string cells exactly (no style/formula left behind).

$ws.Range("D2").Formula = '="60.088.16"'
$ws.Range("E2").Formula = '="  +3.87%  "'

$ws.Range("D3").Formula = '="3.197.10"'
$ws.Range("E3").Formula = '="  +2.39%  "'

$ws.Range("E4").Formula = '="  -0.01%  "'

$ws.Range("D5").Formula = '="537.94"'
$ws.Range("E5").Formula = '="  +0.94%  "'

$ws.Range("D6").Formula = '="145.36"'
$ws.Range("E6").Formula = '="  +5.06%  "'

$ws.Range("D7").Formula = '="0.999"'
$ws.Range("E7").Formula = '="  +0.01%  "'

$ws.Range("E8").Formula = '="  +3.76%  "'

$ws.Range("D9").Formula = '="7.31"'
$ws.Range("E9").Formula = '="  -0.50%  "'

$ws.Range("D10").Formula = '="0.114"'
$ws.Range("E10").Formula = '="  +5.41%  "'

$ws.Range("D11").Formula = '="0.431"'
$ws.Range("E11").Formula = '="  +4.13%  "'

$ws.Range("D12").Formula = '="3.750.73"'
$ws.Range("E12").Formula = '="  +2.57%  "'

$ws.Range("E13").Formula = '="  -0.72%  "'

$ws.Range("D14").Formula = '="0.0000175"'
$ws.Range("E14").Formula = '="  +4.99%  "'

$ws.Range("D15").Formula = '="26.07"'
$ws.Range("E15").Formula = '="  +1.62%  "'

$ws.Range("D16").Formula = '="60.082.73"'
$ws.Range("E16").Formula = '="  +3.68%  "'

$ws.Range("D17").Formula = '="3.208.26"'
$ws.Range("E17").Formula = '="  +2.83%  "'

$ws.Range("D18").Formula = '="6.21"'
$ws.Range("E18").Formula = '="  +0.87%  "'

$ws.Range("D19").Formula = '="13.10"'
$ws.Range("E19").Formula = '="  +2.37%  "'

$ws.Range("B20").Formula = '="BitcoinCash"'
$ws.Range("C20").Formula = '="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"'
$ws.Range("D20").Formula = '="386.05"'
$ws.Range("E20").Formula = '="  +3.05%  "'

$ws.Range("B21").Formula = '="Uniswap"'
$ws.Range("C21").Formula = '="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"'
$ws.Range("D21").Formula = '="8.29"'
$ws.Range("E21").Formula = '="  +2.26%  "'

$ws.Range("D22").Formula = '="1.00"'
$ws.Range("E22").Formula = '="  +0.03%  "'

$ws.Range("D23").Formula = '="0.530"'
$ws.Range("E23").Formula = '="  +4.23%  "'

$ws.Range("D24").Formula = '="70.25"'
$ws.Range("E24").Formula = '="  +0.96%  "'

$ws.Range("B25").Formula = '="Kaspa"'
$ws.Range("C25").Formula = '="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"'
$ws.Range("D25").Formula = '="0.173"'
$ws.Range("E25").Formula = '="  +3.31%  "'

$ws.Range("B26").Formula = '="InternetComputer(DFINITY)"'
$ws.Range("C26").Formula = '="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"'
$ws.Range("D26").Formula = '="8.90"'
$ws.Range("E26").Formula = '="  +15.70%  "'

$ws.Range("D27").Formula = '="0.998"'
$ws.Range("E27").Formula = '="  -0.08%  "'

$ws.Range("D28").Formula = '="0.0₃0909"'
$ws.Range("E28").Formula = '="  +3.00%  "'

$ws.Range("D29").Formula = '="1.92"'
$ws.Range("E29").Formula = '="  +2.66%  "'

$ws.Range("D30").Formula = '="22.44"'
$ws.Range("E30").Formula = '="  +4.01%  "'

$ws.Range("B31").Formula = '="NEARProtocol"'
$ws.Range("C31").Formula = '="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"'
$ws.Range("D31").Formula = '="5.43"'
$ws.Range("E31").Formula = '="  +5.55%  "'

$ws.Range("B32").Formula = '="RenderToken"'
$ws.Range("C32").Formula = '="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"'
$ws.Range("D32").Formula = '="6.16"'
$ws.Range("E32").Formula = '="  +0.48%  "'

$ws.Range("D33").Formula = '="1.21"'
$ws.Range("E33").Formula = '="  +3.04%  "'

$ws.Range("D34").Formula = '="6.48"'
$ws.Range("E34").Formula = '="  +5.05%  "'

$ws.Range("D35").Formula = '="156.30"'
$ws.Range("E35").Formula = '="  -2.83%  "'

$ws.Range("E36").Formula = '="  +2.38%  "'

$ws.Range("B37").Formula = '="EnergySwap"'
$ws.Range("C37").Formula = '="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"'
$ws.Range("D37").Formula = '="25.78"'
$ws.Range("E37").Formula = '="  +1.12%  "'

$ws.Range("B38").Formula = '="Maker"'
$ws.Range("C38").Formula = '="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"'
$ws.Range("D38").Formula = '="2.766.78"'
$ws.Range("E38").Formula = '="  +8.37%  "'

$ws.Range("E39").Formula = '="  +6.25%  "'

$ws.Range("D40").Formula = '="1.69"'
$ws.Range("E40").Formula = '="  +3.47%  "'

$ws.Range("E41").Formula = '="  +2.77%  "'

$ws.Range("D42").Formula = '="0.729"'
$ws.Range("E42").Formula = '="  +4.47%  "'

$ws.Range("D43").Formula = '="39.60"'
$ws.Range("E43").Formula = '="  +2.24%  "'

$ws.Range("D44").Formula = '="0.0288"'
$ws.Range("E44").Formula = '="  +6.73%  "'

$ws.Range("D45").Formula = '="3.243.32"'
$ws.Range("E45").Formula = '="  +2.59%  "'

$ws.Range("D46").Formula = '="1.01"'
$ws.Range("E46").Formula = '="  +3.03%  "'

$ws.Range("D47").Formula = '="6.20"'
$ws.Range("E47").Formula = '="  +0.57%  "'

$ws.Range("E48").Formula = '="  +4.88%  "'

$ws.Range("D49").Formula = '="20.69"'
$ws.Range("E49").Formula = '="  +3.88%  "'

$ws.Range("D50").Formula = '="0.787"'
$ws.Range("E50").Formula = '="  +5.47%  "'

$ws.Range("D51").Formula = '="0.999"'
$ws.Range("E51").Formula = '="  -0.07%  "'

# Collapse every `="literal"` helper formula above into a plain static value
# (text, no formula, no style change) by copying the whole used range onto
# itself with a values-only paste.
$used = $ws.Range("A1:E51")
$used.Copy() | Out-Null
$used.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

